# Insert a new weekly price record for "Choclo" (Macroferia Regional de Talca)
# above the existing row 246, shifting all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 246 (pushes old rows 246..334 down to 247..335)
$ws.Rows.Item(246).Insert()

# Populate the newly inserted row with the new record's data
$ws.Range("A246").Value = 5
$ws.Range("B246").Value = "Macroferia Regional de Talca"
$ws.Range("C246").Value = "Maule"
$ws.Range("D246").Value = 45027
$ws.Range("D246").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E246").Value = 7
$ws.Range("F246").Value = 100112024
$ws.Range("G246").Value = "Choclo"
$ws.Range("H246").Value = "Choclero"
$ws.Range("I246").Value = "Primera"
$ws.Range("J246").Value = 30000
$ws.Range("K246").Value = 350
$ws.Range("L246").Value = 350
$ws.Range("M246").Value = 350
$ws.Range("N246").Value = "$/unidad"
$ws.Range("O246").Value = "Región del Maule"
$ws.Range("P246").Value = 350
$ws.Range("Q246").Value = 1
$ws.Range("R246").Value = "Hortaliza"
